$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 4.2
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("U3").Value = 4.5
$ws.Range("V3").Value = 1.18
$ws.Range("AA3").Value = 5.5
$ws.Range("AH3").Value = 6.5
$ws.Range("AL3").Value = 9.5
$ws.Range("AO3").Value = 41

# Row 4
$ws.Range("G4").Value = 1.75
$ws.Range("J4").Value = 2.4
$ws.Range("AB4").Value = 8
$ws.Range("AG4").Value = 9

# Row 5
$ws.Range("H5").Value = 3.8
$ws.Range("K5").Value = 2.3
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("AK5").Value = 251
